$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1801801801801802
$ws.Range("C2").Value = 0.5990990990990991
$ws.Range("J2").Value = 0.01801801801801802
$ws.Range("P2").Value = 0.1126126126126126
$ws.Range("S2").Value = 0.09009009009009009
$ws.Range("B3").Value = 0.007042253521126761
$ws.Range("C3").Value = 0.06338028169014084
$ws.Range("J3").Value = 0.0352112676056338
$ws.Range("P3").Value = 0.704225352112676
$ws.Range("S3").Value = 0.1901408450704225
$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("O4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.5238095238095238
$ws.Range("S4").Value = 0.3571428571428572
$ws.Range("B6").Value = 0.06145251396648044
$ws.Range("D6").Value = 0.0111731843575419
$ws.Range("F6").Value = 0.05586592178770949
$ws.Range("J6").Value = 0.2290502793296089
$ws.Range("O6").Value = 0.02793296089385475
$ws.Range("Q6").Value = 0.1620111731843575
$ws.Range("R6").Value = 0.0335195530726257
$ws.Range("S6").Value = 0.4189944134078212
$ws.Range("B7").Value = 0.09375
$ws.Range("D7").Value = 0.0125
$ws.Range("F7").Value = 0.05625
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.00625
$ws.Range("Q7").Value = 0.15625
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.4875
$ws.Range("B8").Value = 0.07170542635658915
$ws.Range("D8").Value = 0.01356589147286822
$ws.Range("F8").Value = 0.0562015503875969
$ws.Range("J8").Value = 0.1375968992248062
$ws.Range("O8").Value = 0.02131782945736434
$ws.Range("Q8").Value = 0.1434108527131783
$ws.Range("R8").Value = 0.08914728682170543
$ws.Range("S8").Value = 0.4670542635658915
$ws.Range("B9").Value = 0.1276595744680851
$ws.Range("D9").Value = 0.01063829787234043
$ws.Range("F9").Value = 0.04787234042553191
$ws.Range("J9").Value = 0.1117021276595745
$ws.Range("O9").Value = 0.005319148936170213
$ws.Range("Q9").Value = 0.1648936170212766
$ws.Range("R9").Value = 0.09042553191489362
$ws.Range("S9").Value = 0.4414893617021277
$ws.Range("B10").Value = 0.08348457350272233
$ws.Range("D10").Value = 0.02722323049001815
$ws.Range("F10").Value = 0.06442831215970962
$ws.Range("J10").Value = 0.1152450090744102
$ws.Range("O10").Value = 0.0190562613430127
$ws.Range("Q10").Value = 0.2295825771324864
$ws.Range("R10").Value = 0.07259528130671507
$ws.Range("S10").Value = 0.3883847549909256
$ws.Range("G11").Value = 0.1176470588235294
$ws.Range("J11").Value = 0.08627450980392157
$ws.Range("K11").Value = 0.2352941176470588
$ws.Range("L11").Value = 0.5411764705882353
$ws.Range("S11").Value = 0.0196078431372549
$ws.Range("G12").Value = 0.7342657342657343
$ws.Range("J12").Value = 0.2027972027972028
$ws.Range("K12").Value = 0.01398601398601399
$ws.Range("L12").Value = 0.02797202797202797
$ws.Range("S12").Value = 0.02097902097902098
$ws.Range("G13").Value = 0.825
$ws.Range("J13").Value = 0.15
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.00510204081632653
$ws.Range("H15").Value = 0.2244897959183673
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.2959183673469388
$ws.Range("K15").Value = 0.0663265306122449
$ws.Range("M15").Value = 0.00510204081632653
$ws.Range("O15").Value = 0.06122448979591837
$ws.Range("S15").Value = 0.2704081632653061
$ws.Range("F16").Value = 0.02054794520547945
$ws.Range("H16").Value = 0.2328767123287671
$ws.Range("I16").Value = 0.0958904109589041
$ws.Range("J16").Value = 0.3561643835616438
$ws.Range("K16").Value = 0.08904109589041095
$ws.Range("M16").Value = 0.0273972602739726
$ws.Range("O16").Value = 0.0273972602739726
$ws.Range("S16").Value = 0.1506849315068493
$ws.Range("F17").Value = 0.01970443349753695
$ws.Range("H17").Value = 0.2536945812807882
$ws.Range("I17").Value = 0.07389162561576355
$ws.Range("J17").Value = 0.3940886699507389
$ws.Range("K17").Value = 0.0812807881773399
$ws.Range("M17").Value = 0.01477832512315271
$ws.Range("O17").Value = 0.05911330049261083
$ws.Range("S17").Value = 0.103448275862069
$ws.Range("F18").Value = 0.02515723270440252
$ws.Range("H18").Value = 0.220125786163522
$ws.Range("I18").Value = 0.1132075471698113
$ws.Range("J18").Value = 0.4150943396226415
$ws.Range("K18").Value = 0.0880503144654088
$ws.Range("M18").Value = 0.02515723270440252
$ws.Range("O18").Value = 0.03773584905660377
$ws.Range("S18").Value = 0.07547169811320754
$ws.Range("F19").Value = 0.01045016077170418
$ws.Range("H19").Value = 0.2411575562700965
$ws.Range("I19").Value = 0.09003215434083602
$ws.Range("J19").Value = 0.3512861736334405
$ws.Range("K19").Value = 0.09646302250803858
$ws.Range("M19").Value = 0.02009646302250804
$ws.Range("N19").Value = 0.0008038585209003215
$ws.Range("O19").Value = 0.06591639871382636
$ws.Range("S19").Value = 0.1237942122186495

Write-Host "Updated transition matrix probabilities for added games"
